# Weekly update of "Fruta / hortaliza" data: rows 148-156 get new values
# (data effectively shifts down by one row), row 157 gets new content, and a
# brand-new row 158 is appended (duplicate of the old row 157's tail entry
# but with its own date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 148 ---
$ws.Range("D148").Value2 = 45267
$ws.Range("M148").Value2 = 200
$ws.Range("N148").Value2 = 4000
$ws.Range("O148").Value2 = 5000
$ws.Range("P148").Value2 = 4500
$ws.Range("R148").Value2 = "Región de Ñuble"
$ws.Range("S148").Value2 = 2250

# --- Row 149 ---
$ws.Range("L149").Value2 = "Primera"
$ws.Range("M149").Value2 = 100
$ws.Range("N149").Value2 = 3000
$ws.Range("O149").Value2 = 3500
$ws.Range("P149").Value2 = 3250
$ws.Range("S149").Value2 = 1625

# --- Row 150 ---
$ws.Range("D150").Value2 = 44187
$ws.Range("L150").Value2 = "Segunda"
$ws.Range("M150").Value2 = 50
$ws.Range("N150").Value2 = 2500
$ws.Range("O150").Value2 = 2500
$ws.Range("P150").Value2 = 2500
$ws.Range("R150").Value2 = "Provincia de Curicó"
$ws.Range("S150").Value2 = 1250

# --- Row 151 ---
$ws.Range("D151").Value2 = 44525
$ws.Range("M151").Value2 = 140
$ws.Range("N151").Value2 = 4000
$ws.Range("O151").Value2 = 4500
$ws.Range("P151").Value2 = 4214
$ws.Range("S151").Value2 = 2107

# --- Row 152 ---
$ws.Range("D152").Value2 = 44588
$ws.Range("M152").Value2 = 150
$ws.Range("N152").Value2 = 3000
$ws.Range("O152").Value2 = 3500
$ws.Range("P152").Value2 = 3267
$ws.Range("R152").Value2 = "Provincia de Linares"
$ws.Range("S152").Value2 = 1634

# --- Row 153 ---
$ws.Range("D153").Value2 = 44883
$ws.Range("M153").Value2 = 180
$ws.Range("N153").Value2 = 6000
$ws.Range("O153").Value2 = 6500
$ws.Range("P153").Value2 = 6222
$ws.Range("R153").Value2 = "Región de O'Higgins"
$ws.Range("S153").Value2 = 3111

# --- Row 154 ---
$ws.Range("L154").Value2 = "Primera"
$ws.Range("M154").Value2 = 200
$ws.Range("N154").Value2 = 3500
$ws.Range("O154").Value2 = 4000
$ws.Range("P154").Value2 = 3750
$ws.Range("S154").Value2 = 1875

# --- Row 155 ---
$ws.Range("D155").Value2 = 44558
$ws.Range("L155").Value2 = "Segunda"
$ws.Range("M155").Value2 = 100
$ws.Range("O155").Value2 = 3000
$ws.Range("P155").Value2 = 3000
$ws.Range("S155").Value2 = 1500

# --- Row 156 ---
$ws.Range("D156").Value2 = 44957
$ws.Range("R156").Value2 = "Región de Ñuble"

# --- Row 157 (new content replacing the previous tail row) ---
$ws.Range("L157").Value2 = "Primera"
$ws.Range("M157").Value2 = 200
$ws.Range("N157").Value2 = 3000
$ws.Range("O157").Value2 = 3500
$ws.Range("P157").Value2 = 3250
$ws.Range("S157").Value2 = 1625

# --- Row 158 (brand-new row, appended at the end) ---
$ws.Range("A158").Value2 = 11
$ws.Range("B158").Value2 = "Vega Monumental Concepción"
$ws.Range("C158").Value2 = "Bíobío"
$ws.Range("D158").Value2 = 44897
$ws.Range("D158").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E158").Value2 = 8
$ws.Range("F158").Value2 = "Fruta"
$ws.Range("G158").Value2 = 100101
$ws.Range("H158").Value2 = "Berries"
$ws.Range("I158").Value2 = 100101001
$ws.Range("J158").Value2 = "Arándano (blue)"
$ws.Range("K158").Value2 = "Sin especificar"
$ws.Range("L158").Value2 = "Segunda"
$ws.Range("M158").Value2 = 100
$ws.Range("N158").Value2 = 2800
$ws.Range("O158").Value2 = 2800
$ws.Range("P158").Value2 = 2800
$ws.Range("Q158").Value2 = "$/bandeja 2 kilos"
$ws.Range("R158").Value2 = "Región de O'Higgins"
$ws.Range("S158").Value2 = 1400
$ws.Range("T158").Value2 = 2
